$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - shift header labels
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Data row (row 2)
$ws.Range("C2").Value = "g__Bifidobacterium"
$ws.Range("D2").Value = "g__Bifidobacterium"
$ws.Range("E2").Value = 1
